$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A136").Value = "Ebey's Landing National Historical Reserve"
$ws.Range("A137").Value = "Edgar Allan Poe National Historic Site"
$ws.Range("A138").Value = "Effigy Mounds National Monument"
$ws.Range("A139").Value = "Eisenhower National Historic Site"
$ws.Range("A140").Value = "El Camino Real de los Tejas National Historic Trail"
$ws.Range("A141").Value = "El Camino Real de Tierra Adentro National Historic Trail"
$ws.Range("A142").Value = "El Malpais National Monument"
$ws.Range("A143").Value = "El Morro National Monument"
$ws.Range("A144").Value = "Eleanor Roosevelt National Historic Site"
$ws.Range("A145").Value = "Ellis Island Part of Statue of Liberty National Monument"
$ws.Range("A146").Value = "Erie Canalway National Heritage Corridor"
$ws.Range("A147").Value = "Essex National Heritage Area"
$ws.Range("A148").Value = "Eugene O'Neill National Historic Site"
$ws.Range("A149").Value = "Everglades National Park"
$ws.Range("A151").Value = "Fallen Timbers Battlefield and Fort Miamis National Historic Site"
$ws.Range("A152").Value = "Federal Hall National Memorial"
$ws.Range("A153").Value = "Fire Island National Seashore"
$ws.Range("A154").Value = "First Ladies National Historic Site"
$ws.Range("A155").Value = "First State National Historical Park"
$ws.Range("A156").Value = "Flight 93 National Memorial"
$ws.Range("A157").Value = "Florissant Fossil Beds National Monument"
$ws.Range("A158").Value = "Ford's Theatre"
$ws.Range("A159").Value = "Fort Bowie National Historic Site"
$ws.Range("A160").Value = "Fort Caroline National Memorial"
$ws.Range("A161").Value = "Fort Davis National Historic Site"
$ws.Range("A162").Value = "Fort Donelson National Battlefield"
$ws.Range("A163").Value = "Fort Dupont Park"
$ws.Range("A164").Value = "Fort Foote Park"
$ws.Range("A165").Value = "Fort Frederica National Monument"
$ws.Range("A166").Value = "Fort Laramie National Historic Site"
$ws.Range("A167").Value = "Fort Larned National Historic Site"
$ws.Range("A168").Value = "Fort Matanzas National Monument"
$ws.Range("A169").Value = "Fort McHenry National Monument and Historic Shrine"
$ws.Range("A170").Value = "Fort Monroe National Monument"
$ws.Range("A171").Value = "Fort Necessity National Battlefield"
$ws.Range("A172").Value = "Fort Point National Historic Site"
$ws.Range("A173").Value = "Fort Pulaski National Monument"
$ws.Range("A174").Value = "Fort Raleigh National Historic Site"
$ws.Range("A175").Value = "Fort Scott National Historic Site"
$ws.Range("A176").Value = "Fort Smith National Historic Site"
$ws.Range("A177").Value = "Fort Stanwix National Monument"
$ws.Range("A178").Value = "Fort Sumter and Fort Moultrie National Historical Park"
$ws.Range("A179").Value = "Fort Union National Monument"
$ws.Range("A180").Value = "Fort Union Trading Post National Historic Site"
$ws.Range("A181").Value = "Fort Vancouver National Historic Site"
$ws.Range("A182").Value = "Fort Washington Park"
$ws.Range("A183").Value = "Fossil Butte National Monument"
$ws.Range("A184").Value = "Franklin Delano Roosevelt Memorial"
$ws.Range("A185").Value = "Frederick Douglass National Historic Site"
$ws.Range("A186").Value = "Frederick Law Olmsted National Historic Site"
$ws.Range("A187").Value = "Fredericksburg & Spotsylvania National Military Park"
$ws.Range("A188").Value = "Freedom Riders National Monument"
$ws.Range("A189").Value = "Freedom's Way National Heritage Area"
$ws.Range("A190").Value = "Friendship Hill National Historic Site"

try {
    [void]$ws.Range("A190").Select()
    $excel.ActiveWindow.ScrollRow = 179
} catch {
}
